$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9,1).Value2 = "J1"
$ws.Cells.Item(9,2).Value2 = "-"
$ws.Cells.Item(9,3).Value2 = "-"
$ws.Cells.Item(9,4).Value2 = "-"
$ws.Cells.Item(9,5).Value2 = "Micro-Fit connector"
$ws.Cells.Item(9,6).Value2 = "Molex"
$ws.Cells.Item(9,7).Value2 = "43045-0210"
$ws.Cells.Item(9,8).Value2 = 3103032
